$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "51.553.07"
Set-TextValue $ws "E2" "  -1.46%  "
Set-TextValue $ws "D3" "2.952.54"
Set-TextValue $ws "E3" "  -1.97%  "
Set-TextValue $ws "D4" "0.998"
Set-TextValue $ws "E4" "  -0.22%  "
Set-TextValue $ws "D5" "377.18"
Set-TextValue $ws "E5" "  +6.06%  "
Set-TextValue $ws "D6" "104.82"
Set-TextValue $ws "E6" "  -3.46%  "
Set-TextValue $ws "D7" "0.546"
Set-TextValue $ws "E7" "  -3.28%  "
Set-TextValue $ws "D8" "0.998"
Set-TextValue $ws "E8" "  -0.20%  "
Set-TextValue $ws "E9" "  -4.15%  "
Set-TextValue $ws "D10" "37.45"
Set-TextValue $ws "E10" "  -2.84%  "
Set-TextValue $ws "D11" "0.140"
Set-TextValue $ws "E11" "  +0.26%  "
Set-TextValue $ws "D12" "0.0843"
Set-TextValue $ws "E12" "  -2.11%  "
Set-TextValue $ws "D13" "18.48"
Set-TextValue $ws "E13" "  -4.53%  "
Set-TextValue $ws "D14" "3.412.87"
Set-TextValue $ws "E14" "  -2.15%  "
Set-TextValue $ws "D15" "7.43"
Set-TextValue $ws "E15" "  -4.37%  "
Set-TextValue $ws "D16" "2.946.20"
Set-TextValue $ws "E16" "  -1.92%  "
Set-TextValue $ws "D17" "0.945"
Set-TextValue $ws "E17" "  -7.89%  "
Set-TextValue $ws "D18" "51.542.17"
Set-TextValue $ws "E18" "  -1.57%  "
Set-TextValue $ws "E19" "  -6.00%  "
Set-TextValue $ws "D20" "7.36"
Set-TextValue $ws "E20" "  -2.56%  "
Set-TextValue $ws "D21" "13.11"
Set-TextValue $ws "E21" "  -4.52%  "
Set-TextValue $ws "D22" "0.0₃0951"
Set-TextValue $ws "E22" "  -2.70%  "
Set-TextValue $ws "D23" "68.83"
Set-TextValue $ws "E23" "  -1.26%  "
Set-TextValue $ws "D24" "262.46"
Set-TextValue $ws "E24" "  -1.22%  "
Set-TextValue $ws "D25" "2.73"
Set-TextValue $ws "E25" "  -0.97%  "
Set-TextValue $ws "D26" "0.171"
Set-TextValue $ws "E26" "  -4.67%  "
Set-TextValue $ws "B27" "LEO"
Set-TextValue $ws "C27" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws "D27" "4.16"
Set-TextValue $ws "E27" "  -2.93%  "
Set-TextValue $ws "B28" "Dai"
Set-TextValue $ws "C28" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws "D28" "1.00"
Set-TextValue $ws "E28" "  +0.06%  "
Set-TextValue $ws "B29" "EthereumClassic"
Set-TextValue $ws "C29" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws "D29" "25.99"
Set-TextValue $ws "E29" "  -4.04%  "
Set-TextValue $ws "D30" "6.96"
Set-TextValue $ws "E30" "  +7.37%  "
Set-TextValue $ws "B31" "Filecoin"
Set-TextValue $ws "C31" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws "D31" "7.18"
Set-TextValue $ws "E31" "  -6.24%  "
Set-TextValue $ws "B32" "Hedera"
Set-TextValue $ws "C32" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws "D32" "0.104"
Set-TextValue $ws "E32" "  -3.18%  "
Set-TextValue $ws "B33" "Cosmos"
Set-TextValue $ws "C33" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws "D33" "9.94"
Set-TextValue $ws "E33" "  -4.14%  "
Set-TextValue $ws "B34" "InjectiveProtocol"
Set-TextValue $ws "C34" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws "D34" "34.91"
Set-TextValue $ws "E34" "  -4.69%  "
Set-TextValue $ws "B35" "Toncoin"
Set-TextValue $ws "C35" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws "D35" "2.12"
Set-TextValue $ws "E35" "  -3.55%  "
Set-TextValue $ws "B36" "OKB"
Set-TextValue $ws "C36" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws "D36" "50.23"
Set-TextValue $ws "E36" "  -1.47%  "
Set-TextValue $ws "B37" "VeChain"
Set-TextValue $ws "C37" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D37" "0.0430"
Set-TextValue $ws "E37" "  -2.92%  "
Set-TextValue $ws "B38" "FirstDigitalUSD"
Set-TextValue $ws "C38" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws "D38" "1.00"
Set-TextValue $ws "E38" "  +0.53%  "
Set-TextValue $ws "B39" "LidoDAOToken"
Set-TextValue $ws "C39" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws "D39" "3.05"
Set-TextValue $ws "E39" "  -6.20%  "
Set-TextValue $ws "B40" "Celestia"
Set-TextValue $ws "C40" "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws "D40" "17.26"
Set-TextValue $ws "E40" "  -4.25%  "
Set-TextValue $ws "B41" "Stacks"
Set-TextValue $ws "C41" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws "D41" "2.62"
Set-TextValue $ws "E41" "  -3.42%  "
Set-TextValue $ws "B42" "ARBITRUM"
Set-TextValue $ws "C42" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws "D42" "1.87"
Set-TextValue $ws "E42" "  -7.29%  "
Set-TextValue $ws "B43" "Stellar"
Set-TextValue $ws "C43" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws "D43" "0.114"
Set-TextValue $ws "E43" "  -2.85%  "
Set-TextValue $ws "B44" "EnergySwap"
Set-TextValue $ws "C44" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws "D44" "22.21"
Set-TextValue $ws "E44" "  -3.42%  "
Set-TextValue $ws "B45" "Monero"
Set-TextValue $ws "C45" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws "D45" "119.55"
Set-TextValue $ws "E45" "  -3.10%  "
Set-TextValue $ws "B46" "WEMIXToken"
Set-TextValue $ws "C46" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws "D46" "2.11"
Set-TextValue $ws "E46" "  -2.31%  "
Set-TextValue $ws "B47" "Maker"
Set-TextValue $ws "C47" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws "D47" "2.041.35"
Set-TextValue $ws "E47" "  -4.41%  "
Set-TextValue $ws "B48" "ApeXProtocol"
Set-TextValue $ws "C48" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws "D48" "2.31"
Set-TextValue $ws "E48" "  -3.77%  "
Set-TextValue $ws "E49" "  +6.67%  "
Set-TextValue $ws "B50" "NEARProtocol"
Set-TextValue $ws "C50" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws "D50" "3.22"
Set-TextValue $ws "E50" "  -5.47%  "
Set-TextValue $ws "B51" "RocketPoolETH"
Set-TextValue $ws "C51" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue $ws "D51" "3.238.58"
Set-TextValue $ws "E51" "  -2.25%  "
